# add area to Q files stn6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2: first area segment + totals
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3: area segment using previous depth
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15: area segments (fill down the same relative formula)
for ($r = 4; $r -le 15; $r++) {
    $prev = $r - 1
    $ws.Range("G$r").Formula = "=(D$r-D$prev)*B$r/100"
}

# Match the author's final selection (the new Atotal/Qtotal summary cells)
$ws.Range("J2:K2").Select()

$wb.Save()
